$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 and 3 (delete higher index first so row numbers don't shift)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Row 2 updates
$ws.Cells.Item(2, 1).Value = 99

$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "0008-08-08"
$ws.Cells.Item(2, 2).Style = "Normal"

$ws.Cells.Item(2, 3).Value = "t de 2"

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "8"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).Value = $true

# New column F - header "Status" styled like the other header cells
$ws.Cells.Item(1, 6).Value = "Status"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Cells.Item(2, 6).Value = "Expedição"
